$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 676.3333
$ws.Range("I2").Value = 486
$ws.Range("J2").Value = 866.6667
$ws.Range("K2").Value = 486
$ws.Range("L2").Value = 866.6667
$ws.Range("M2").Value = -373
$ws.Range("N2").Value = -1092.6667
$ws.Range("H5").Value = 319
$ws.Range("I5").Value = 111.6
$ws.Range("J5").Value = 1010.3333
$ws.Range("K5").Value = 111.6
$ws.Range("L5").Value = 1010.3333
$ws.Range("M5").Value = 3.400000000000006
$ws.Range("N5").Value = -1240.3333
$ws.Range("H33").Value = 203.5
$ws.Range("I33").Value = 207.33333
$ws.Range("J33").Value = 192
$ws.Range("K33").Value = 207.33333
$ws.Range("L33").Value = 192
$ws.Range("M33").Value = 21.66667000000001
$ws.Range("N33").Value = -650
$ws.Range("H51").Value = 24499.5
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H70").Value = 1933
$ws.Range("J70").Value = 1933
$ws.Range("L70").Value = 5799
$ws.Range("N70").Value = -6339
$ws.Range("H73").Value = 1933
$ws.Range("J73").Value = 1933
$ws.Range("L73").Value = 5799
$ws.Range("N73").Value = -7671
$ws.Range("H101").Value = 655.75
$ws.Range("I101").Value = 657
$ws.Range("J101").Value = 652
$ws.Range("K101").Value = 1971
$ws.Range("L101").Value = 1956
$ws.Range("M101").Value = -349
$ws.Range("N101").Value = -5200
$ws.Range("H138").Value = 6238.079
$ws.Range("I138").Value = 1498.0667
$ws.Range("J138").Value = 9329.392
$ws.Range("K138").Value = 4494.2001
$ws.Range("L138").Value = 27988.176
$ws.Range("M138").Value = 645.7999
$ws.Range("N138").Value = -38268.176

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4682.4165
$ws.Range("I61").Value = 4835.636
$ws.Range("K61").Value = 4835.636
$ws.Range("M61").Value = -4623.636
$ws.Range("H74").Value = 1278.7
$ws.Range("I74").Value = 1426.7142
$ws.Range("J74").Value = 933.3333
$ws.Range("K74").Value = 1426.7142
$ws.Range("L74").Value = 933.3333
$ws.Range("M74").Value = -552.7141999999999
$ws.Range("N74").Value = -2681.3333
$ws.Range("H77").Value = 1278.7
$ws.Range("I77").Value = 1426.7142
$ws.Range("J77").Value = 933.3333
$ws.Range("K77").Value = 7133.571
$ws.Range("L77").Value = 4666.6665
$ws.Range("M77").Value = -2765.571
$ws.Range("N77").Value = -13402.6665
$ws.Range("H132").Value = 2256.5652
$ws.Range("I132").Value = 1860.6842
$ws.Range("K132").Value = 5582.0526
$ws.Range("M132").Value = -3052.0526
$ws.Range("H136").Value = 4682.4165
$ws.Range("I136").Value = 4835.636
$ws.Range("K136").Value = 14506.908
$ws.Range("M136").Value = -11956.908

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 328.625
$ws.Range("J7").Value = 484.2
$ws.Range("L7").Value = 484.2
$ws.Range("N7").Value = -710.2
$ws.Range("H22").Value = 546.8182
$ws.Range("I22").Value = 315.8
$ws.Range("K22").Value = 315.8
$ws.Range("M22").Value = 34.19999999999999
$ws.Range("H62").Value = 4153.7
$ws.Range("I62").Value = 3791.8572
$ws.Range("K62").Value = 3791.8572
$ws.Range("M62").Value = -3167.8572
$ws.Range("H65").Value = 4153.7
$ws.Range("I65").Value = 3791.8572
$ws.Range("K65").Value = 18959.286
$ws.Range("M65").Value = -15839.286
$ws.Range("H134").Value = 3823.2144
$ws.Range("I134").Value = 3794
$ws.Range("K134").Value = 11382
$ws.Range("M134").Value = -8847

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 900
$ws.Range("J5").Value = 900
$ws.Range("L5").Value = 2700
$ws.Range("N5").Value = -2924
$ws.Range("H113").Value = 2995.6667
$ws.Range("J113").Value = 2995.6667
$ws.Range("L113").Value = 8987.000100000001
$ws.Range("N113").Value = -13327.0001
$ws.Range("H125").Value = 7021.696
$ws.Range("J125").Value = 7500
$ws.Range("L125").Value = 22500
$ws.Range("N125").Value = -32340
$ws.Range("H135").Value = 900
$ws.Range("J135").Value = 900
$ws.Range("L135").Value = 8100
$ws.Range("N135").Value = -13170

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 929.1667
$ws.Range("I122").Value = 919
$ws.Range("J122").Value = 949.5
$ws.Range("K122").Value = 2757
$ws.Range("L122").Value = 2848.5
$ws.Range("M122").Value = -307
$ws.Range("N122").Value = -7748.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 5000
$ws.Range("J38").Value = 5000
$ws.Range("L38").Value = 5000
$ws.Range("N38").Value = -5820
$ws.Range("H42").Value = 34999
$ws.Range("J42").Value = 34999
$ws.Range("L42").Value = 34999
$ws.Range("N42").Value = -36125
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H49").Value = 34999
$ws.Range("J49").Value = 34999
$ws.Range("L49").Value = 34999
$ws.Range("N49").Value = -35293
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H132").Value = 3131.7222
$ws.Range("I132").Value = 1938.7
$ws.Range("K132").Value = 5816.1
$ws.Range("M132").Value = -3286.1
$ws.Range("H136").Value = 4223.778
$ws.Range("I136").Value = 4223.778
$ws.Range("K136").Value = 12671.334
$ws.Range("M136").Value = -10121.334

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2641
$ws.Range("I132").Value = 1748.75
$ws.Range("J132").Value = 3938.818
$ws.Range("K132").Value = 5246.25
$ws.Range("L132").Value = 11816.454
$ws.Range("M132").Value = -2716.25
$ws.Range("N132").Value = -16876.454
